$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.241999999999998
$ws.Range("C3").Value = -10.30289999999999
$ws.Range("D3").Value = -6.774199999999991
$ws.Range("C4").Value = -12.409
$ws.Range("E8").Value = 16.0583
$ws.Range("D9").Value = -7.208799999999997
$ws.Range("A11").Value = -21.64249999999999
$ws.Range("E11").Value = 16.5879
$ws.Range("A12").Value = -21.56519999999999
$ws.Range("C14").Value = -13.389
$ws.Range("E14").Value = 16.75730000000001
$ws.Range("A15").Value = -21.54189999999999
$ws.Range("D15").Value = -8.072999999999997
$ws.Range("E15").Value = 16.0691
$ws.Range("E17").Value = 16.88630000000001
$ws.Range("D19").Value = -7.314199999999997
$ws.Range("D20").Value = -7.440100000000001
$ws.Range("D25").Value = -7.540100000000002
$ws.Range("C26").Value = -12.97660000000001
$ws.Range("E26").Value = 15.70029999999999
$ws.Range("A27").Value = -21.8837
$ws.Range("D27").Value = -8.692299999999998
$ws.Range("A28").Value = -21.70759999999999
$ws.Range("D28").Value = -7.839900000000003
$ws.Range("D30").Value = -7.116900000000005
$ws.Range("A31").Value = -21.43780000000001
$ws.Range("C31").Value = -13.0004
$ws.Range("A32").Value = -21.6377
$ws.Range("D32").Value = -8.333000000000006
$ws.Range("C35").Value = -12.4953
$ws.Range("A36").Value = -20.826
$ws.Range("E36").Value = 15.8336
$ws.Range("C37").Value = -14.09209999999999
$ws.Range("A38").Value = -19.4982
$ws.Range("C39").Value = -12.88810000000001
$ws.Range("C40").Value = -13.99970000000001
$ws.Range("E42").Value = 16.64979999999999
$ws.Range("D44").Value = -7.217300000000002
$ws.Range("C45").Value = -13.98329999999999
$ws.Range("A46").Value = -21.4558
$ws.Range("D47").Value = -7.577200000000001
$ws.Range("C52").Value = -10.9131
$ws.Range("A54").Value = -21.4922
$ws.Range("A55").Value = -22.4062
$ws.Range("A56").Value = -22.21490000000001
$ws.Range("C57").Value = -14.91969999999998
$ws.Range("D58").Value = -8.020500000000004
$ws.Range("D62").Value = -8.296099999999999
$ws.Range("E64").Value = 17.4993
$ws.Range("A67").Value = -21.56349999999998
$ws.Range("E68").Value = 16.84940000000001
$ws.Range("A69").Value = -21.71969999999997
$ws.Range("A72").Value = -21.50189999999998
$ws.Range("A73").Value = -19.97290000000001
$ws.Range("D77").Value = -5.5103
$ws.Range("D78").Value = -7.506500000000002
$ws.Range("E79").Value = 18.19090000000002
$ws.Range("C81").Value = -13.05599999999999
$ws.Range("A83").Value = -21.60769999999999
$ws.Range("C83").Value = -12.2716
$ws.Range("D84").Value = -8.673199999999996
$ws.Range("A86").Value = -22.27650000000002
$ws.Range("D89").Value = -6.000199999999999
$ws.Range("E89").Value = 18.82190000000002
$ws.Range("A91").Value = -21.47500000000001
$ws.Range("D91").Value = -5.984799999999999
$ws.Range("D92").Value = -5.947299999999997
$ws.Range("A93").Value = -21.24329999999999
$ws.Range("D96").Value = -7.276900000000005
$ws.Range("A99").Value = -20.55919999999999
$ws.Range("C100").Value = -12.95159999999999
$ws.Range("C102").Value = -14.54239999999999
$ws.Range("D102").Value = -7.8569
